$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1807.7858
$ws.Range("J19").Value = 3198.1667
$ws.Range("L19").Value = 3198.1667
$ws.Range("N19").Value = -3548.1667

# Hunk 1: ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2480
$ws.Range("I40").Value = 2460.1667
$ws.Range("J40").Value = 2499.8333
$ws.Range("K40").Value = 2460.1667
$ws.Range("L40").Value = 2499.8333
$ws.Range("M40").Value = -2285.1667
$ws.Range("N40").Value = -2849.8333

# Hunk 2: ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1383.7333
$ws.Range("I43").Value = 965.5
$ws.Range("J43").Value = 1535.8182
$ws.Range("K43").Value = 965.5
$ws.Range("L43").Value = 1535.8182
$ws.Range("M43").Value = -896.5
$ws.Range("N43").Value = -1673.8182

# Hunk 3: ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2346381.8
$ws.Range("I76").Value = 3349003.2
$ws.Range("K76").Value = 3349003.2
$ws.Range("M76").Value = -3348688.2

# Hunk 4: ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2346381.8
$ws.Range("I79").Value = 3349003.2
$ws.Range("K79").Value = 3349003.2
$ws.Range("M79").Value = -3347911.2

# Hunk 5: ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4718.5713
$ws.Range("J112").Value = 5305
$ws.Range("L112").Value = 15915
$ws.Range("N112").Value = -18131

# Hunk 6: ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 999
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -6491

# Hunk 7: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1159.037
$ws.Range("I132").Value = 1172.8462
$ws.Range("K132").Value = 3518.5386
$ws.Range("M132").Value = -988.5385999999999

# Hunk 8: ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 47616.668
$ws.Range("J139").Value = 47616.668
$ws.Range("L139").Value = 47616.668
$ws.Range("N139").Value = -57896.668

# Hunk 9: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5438.8716
$ws.Range("I32").Value = 4128.68
$ws.Range("J32").Value = 7778.5
$ws.Range("K32").Value = 4128.68
$ws.Range("L32").Value = 7778.5
$ws.Range("M32").Value = -3841.68
$ws.Range("N32").Value = -8352.5

# Hunk 10: ARM row 60
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 13500

# Hunk 11: ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 881.8333
$ws.Range("I97").Value = 458.2
$ws.Range("K97").Value = 458.2
$ws.Range("M97").Value = 37.80000000000001

# Hunk 12: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1892.909
$ws.Range("I132").Value = 1202.6666
$ws.Range("K132").Value = 3607.9998
$ws.Range("M132").Value = -1077.9998

# Hunk 13: BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2751
$ws.Range("I7").Value = 2751
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2751
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2638
$ws.Range("N7").ClearContents()

# Hunk 14: BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2007.1428
$ws.Range("I20").Value = 1751.1
$ws.Range("K20").Value = 1751.1
$ws.Range("M20").Value = -1504.1

# Hunk 15: BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 609.2857
$ws.Range("I94").Value = 540.2
$ws.Range("J94").Value = 782
$ws.Range("K94").Value = 540.2
$ws.Range("L94").Value = 782
$ws.Range("M94").Value = -89.20000000000005
$ws.Range("N94").Value = -1684

# Hunk 16: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3250
$ws.Range("I107").Value = 3250
$ws.Range("K107").Value = 3250
$ws.Range("M107").Value = -1330

# Hunk 17: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6729.4346
$ws.Range("I134").Value = 7147.514
$ws.Range("J134").Value = 5399.1816
$ws.Range("K134").Value = 21442.542
$ws.Range("L134").Value = 16197.5448
$ws.Range("M134").Value = -18907.542
$ws.Range("N134").Value = -21267.5448

# Hunk 18: CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90.833336
$ws.Range("I7").Value = 114.25
$ws.Range("J7").Value = 44
$ws.Range("K7").Value = 114.25
$ws.Range("L7").Value = 44
$ws.Range("M7").Value = -1.25
$ws.Range("N7").Value = -270

# Hunk 19: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2462.5334
$ws.Range("I31").Value = 1180.9565
$ws.Range("K31").Value = 1180.9565
$ws.Range("M31").Value = -885.9565

# Hunk 20: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2462.5334
$ws.Range("I34").Value = 1180.9565
$ws.Range("K34").Value = 1180.9565
$ws.Range("M34").Value = -978.9565

# Hunk 21: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1080.5
$ws.Range("I58").Value = 918.1
$ws.Range("J58").Value = 1486.5
$ws.Range("K58").Value = 918.1
$ws.Range("L58").Value = 1486.5
$ws.Range("M58").Value = -715.1
$ws.Range("N58").Value = -1892.5

# Hunk 22: CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1875.6
$ws.Range("I122").Value = 2385.8572
$ws.Range("J122").Value = 1429.125
$ws.Range("K122").Value = 7157.571599999999
$ws.Range("L122").Value = 4287.375
$ws.Range("M122").Value = -4707.571599999999
$ws.Range("N122").Value = -9187.375

# Hunk 23: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1080.5
$ws.Range("I136").Value = 918.1
$ws.Range("J136").Value = 1486.5
$ws.Range("K136").Value = 2754.3
$ws.Range("L136").Value = 4459.5
$ws.Range("M136").Value = -204.3000000000002
$ws.Range("N136").Value = -9559.5

# Hunk 24: CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 165

# Hunk 25: CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9232.833000000001
$ws.Range("I113").Value = 50352
$ws.Range("J113").Value = 1009
$ws.Range("K113").Value = 151056
$ws.Range("L113").Value = 3027
$ws.Range("M113").Value = -148886
$ws.Range("N113").Value = -7367

# Hunk 26: CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 773.6
$ws.Range("I131").Value = 337.625
$ws.Range("J131").Value = 811.51086
$ws.Range("K131").Value = 1012.875
$ws.Range("L131").Value = 2434.53258
$ws.Range("M131").Value = 4027.125
$ws.Range("N131").Value = -12514.53258

# Hunk 27: GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5322.1113
$ws.Range("J70").Value = 4300
$ws.Range("L70").Value = 4300
$ws.Range("N70").Value = -4840

# Hunk 28: GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5322.1113
$ws.Range("J73").Value = 4300
$ws.Range("L73").Value = 4300
$ws.Range("N73").Value = -6172

# Hunk 29: GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2428.4583
$ws.Range("I102").Value = 2762.7273
$ws.Range("K102").Value = 2762.7273
$ws.Range("M102").Value = -1140.7273

# Hunk 30: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2295.2964
$ws.Range("I122").Value = 2233.923
$ws.Range("J122").Value = 2352.2856
$ws.Range("K122").Value = 6701.768999999999
$ws.Range("L122").Value = 7056.8568
$ws.Range("M122").Value = -4251.768999999999
$ws.Range("N122").Value = -11956.8568

# Hunk 31: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4917.839
$ws.Range("I132").Value = 4025.4583
$ws.Range("K132").Value = 12076.3749
$ws.Range("M132").Value = -9546.374899999999

# Hunk 32: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7768
$ws.Range("I7").Value = 10004
$ws.Range("K7").Value = 10004
$ws.Range("M7").Value = -9892

# Hunk 33: LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705

# Hunk 34: LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893

# Hunk 35: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5515
$ws.Range("I122").Value = 3875
$ws.Range("J122").Value = 6335
$ws.Range("K122").Value = 11625
$ws.Range("L122").Value = 19005
$ws.Range("M122").Value = -9175
$ws.Range("N122").Value = -23905

# Hunk 36: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7768
$ws.Range("I126").Value = 10004
$ws.Range("K126").Value = 30012
$ws.Range("M126").Value = -27542

# Hunk 37: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 117768.91
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

Write-Host "All 38 hunks applied"
